# Apply "I0 and IF added" edit: add columns I (I0) and J (IF) to the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new header cells I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered alignment) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-68): fill in the I0 and J0 values for each row ---
$iValues = @(8,8,9,8,8,6,8,7,6,5,9,8,8,6,5,7,8,7,6,7,7,5,8,7,7,6,1,7,7,4,3,7,7,1,1,6,8,7,6,5,6,5,7,7,1,1,5,1,7,8,8,9,7,7,6,7,5,6,5,8,7,6,1,1,4,1,1)
$jValues = @(8,8,9,8,8,6,8,9,7,6,10,8,8,7,7,7,9,7,7,8,7,7,9,8,7,6,4,7,7,6,5,8,8,6,4,6,8,7,6,6,6,5,8,7,4,4,5,4,7,8,11,9,8,8,6,9,6,6,6,8,7,6,4,4,6,3,2)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}

Write-Output "applied I0/IF columns"
